# Updated cryptos list values (price + volume columns, and a block of
# reordered coin rows 34-38) to match the refreshed source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the literal string into the cell as TEXT (matching the workbook's
    # original inline-string cells) instead of letting Excel auto-coerce
    # numeric-looking strings (e.g. "295.39") into a floating point number.
    # Resetting the style back to "Normal" afterwards avoids leaving a stray
    # text-numberformat style behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "42.417.63"
Set-TextValue $ws.Range("E2") "  -8.57%  "
Set-TextValue $ws.Range("D3") "2.503.45"
Set-TextValue $ws.Range("E3") "  -4.18%  "
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.15%  "
Set-TextValue $ws.Range("D5") "295.39"
Set-TextValue $ws.Range("E5") "  -3.62%  "
Set-TextValue $ws.Range("D6") "92.47"
Set-TextValue $ws.Range("E6") "  -7.70%  "
Set-TextValue $ws.Range("E7") "  -5.71%  "
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  +0.08%  "
Set-TextValue $ws.Range("D9") "0.544"
Set-TextValue $ws.Range("E9") "  -6.16%  "
Set-TextValue $ws.Range("D10") "35.98"
Set-TextValue $ws.Range("E10") "  -8.70%  "
Set-TextValue $ws.Range("E11") "  -5.58%  "
Set-TextValue $ws.Range("E12") "  -6.41%  "
Set-TextValue $ws.Range("E13") "  +0.13%  "
Set-TextValue $ws.Range("D14") "2.890.10"
Set-TextValue $ws.Range("E14") "  -3.74%  "
Set-TextValue $ws.Range("D15") "2.497.00"
Set-TextValue $ws.Range("E15") "  -4.21%  "
Set-TextValue $ws.Range("E16") "  -6.60%  "
Set-TextValue $ws.Range("E17") "  -6.67%  "
Set-TextValue $ws.Range("D18") "42.403.09"
Set-TextValue $ws.Range("E18") "  -8.74%  "
Set-TextValue $ws.Range("D19") "0.0₃0956"
Set-TextValue $ws.Range("E19") "  -5.33%  "
Set-TextValue $ws.Range("D20") "6.47"
Set-TextValue $ws.Range("E20") "  -3.81%  "
Set-TextValue $ws.Range("D21") "12.21"
Set-TextValue $ws.Range("E21") "  -5.86%  "
Set-TextValue $ws.Range("D22") "71.96"
Set-TextValue $ws.Range("E22") "  +0.53%  "
Set-TextValue $ws.Range("D23") "256.80"
Set-TextValue $ws.Range("E23") "  -5.94%  "
Set-TextValue $ws.Range("D24") "2.87"
Set-TextValue $ws.Range("E24") "  -5.34%  "
Set-TextValue $ws.Range("D25") "2.10"
Set-TextValue $ws.Range("E25") "  -3.49%  "
Set-TextValue $ws.Range("E26") "  -2.79%  "
Set-TextValue $ws.Range("E27") "  +0.20%  "
Set-TextValue $ws.Range("D28") "9.86"
Set-TextValue $ws.Range("E28") "  -6.88%  "
Set-TextValue $ws.Range("D29") "2.12"
Set-TextValue $ws.Range("E29") "  -4.26%  "
Set-TextValue $ws.Range("D30") "36.60"
Set-TextValue $ws.Range("E30") "  -4.99%  "
Set-TextValue $ws.Range("D31") "5.96"
Set-TextValue $ws.Range("E31") "  -5.97%  "
Set-TextValue $ws.Range("D32") "3.45"
Set-TextValue $ws.Range("E32") "  -5.74%  "
Set-TextValue $ws.Range("D33") "150.75"
Set-TextValue $ws.Range("E33") "  -0.77%  "
Set-TextValue $ws.Range("B34") "ARBITRUM"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D34") "2.17"
Set-TextValue $ws.Range("E34") "  -2.98%  "
Set-TextValue $ws.Range("B35") "WEMIXToken"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D35") "2.70"
Set-TextValue $ws.Range("E35") "  -5.52%  "
Set-TextValue $ws.Range("B36") "Hedera"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.0791"
Set-TextValue $ws.Range("E36") "  -5.45%  "
Set-TextValue $ws.Range("B37") "Kaspa"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D37") "0.113"
Set-TextValue $ws.Range("E37") "  -7.89%  "
Set-TextValue $ws.Range("B38") "EnergySwap"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D38") "24.29"
Set-TextValue $ws.Range("E38") "  +4.19%  "
Set-TextValue $ws.Range("E39") "  -4.24%  "
Set-TextValue $ws.Range("D40") "16.34"
Set-TextValue $ws.Range("E40") "  +2.89%  "
Set-TextValue $ws.Range("E41") "  -6.00%  "
Set-TextValue $ws.Range("D42") "0.0307"
Set-TextValue $ws.Range("E42") "  -7.03%  "
Set-TextValue $ws.Range("D43") "3.80"
Set-TextValue $ws.Range("E43") "  -6.63%  "
Set-TextValue $ws.Range("D44") "2.000.74"
Set-TextValue $ws.Range("E44") "  -6.01%  "
Set-TextValue $ws.Range("D45") "0.997"
Set-TextValue $ws.Range("E45") "  -0.11%  "
Set-TextValue $ws.Range("D46") "85.14"
Set-TextValue $ws.Range("E46") "  -9.53%  "
Set-TextValue $ws.Range("D47") "1.61"
Set-TextValue $ws.Range("E47") "  +1.67%  "
Set-TextValue $ws.Range("D48") "8.77"
Set-TextValue $ws.Range("E48") "  -8.24%  "
Set-TextValue $ws.Range("D49") "2.743.73"
Set-TextValue $ws.Range("E49") "  -4.05%  "
Set-TextValue $ws.Range("D50") "102.21"
Set-TextValue $ws.Range("E50") "  -6.15%  "
Set-TextValue $ws.Range("D51") "1.64"
Set-TextValue $ws.Range("E51") "  -7.95%  "
